# Update odds values on Sheet1 to reflect the latest Betfair Back/Lay data
# for Jogos_do_Dia_Betfair_Back_Lay_2026-01-19.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L2").Value = 1.3
$ws.Range("F3").Value = 2.68
$ws.Range("G3").Value = 2.82
$ws.Range("K3").Value = 3.95
$ws.Range("Q3").Value = 1.71
$ws.Range("S3").Value = 2.8
$ws.Range("U3").Value = 2.34
$ws.Range("W3").Value = 1.55
$ws.Range("AA4").Value = 250
$ws.Range("F4").Value = 1.54
$ws.Range("I4").Value = 7.4
$ws.Range("AH5").Value = 29
$ws.Range("H5").Value = 1.83
$ws.Range("I5").Value = 1.89
$ws.Range("K5").Value = 4
$ws.Range("S5").Value = 4.1
$ws.Range("U5").Value = 1.85
$ws.Range("V5").Value = 2.1
$ws.Range("Z5").Value = 10.5
$ws.Range("AC6").Value = 13
$ws.Range("G6").Value = 1.51
$ws.Range("K6").Value = 5.3
$ws.Range("L6").Value = 1.36
$ws.Range("O6").Value = 1.37
$ws.Range("W6").Value = 2.92
$ws.Range("X6").Value = 1000
$ws.Range("G7").Value = 1.92
$ws.Range("H7").Value = 4.8
$ws.Range("I7").Value = 6.8
$ws.Range("L7").Value = 1.35
$ws.Range("Q7").Value = 1.89
$ws.Range("S7").Value = 3.55
$ws.Range("V7").Value = 1.19
$ws.Range("W7").Value = 2.08
$ws.Range("X7").Value = 90
$ws.Range("O8").Value = 1.38
$ws.Range("V8").Value = 1.26
$ws.Range("AG9").Value = 10
$ws.Range("N9").Value = 6.4
$ws.Range("P9").Value = 2.88
$ws.Range("J10").Value = 3.1
$ws.Range("M10").Value = 1.05
$ws.Range("S10").Value = 2.78
$ws.Range("G11").Value = 2.72
$ws.Range("T11").Value = 2.02
$ws.Range("N12").Value = 3.3
$ws.Range("O12").Value = 1.35
$ws.Range("P12").Value = 1.73
$ws.Range("Q12").Value = 2.12
$ws.Range("R12").Value = 1.27
$ws.Range("T12").Value = 1.9
$ws.Range("U12").Value = 1.89
$ws.Range("AC13").Value = 8.199999999999999
$ws.Range("AD13").Value = 1000
$ws.Range("AF13").Value = 1000
$ws.Range("AO13").Value = 75
$ws.Range("F13").Value = 2.04
$ws.Range("J13").Value = 3.35
$ws.Range("K13").Value = 3.8
$ws.Range("W13").Value = 1.84
$ws.Range("Y13").Value = 1000
$ws.Range("J14").Value = 3.25
$ws.Range("S14").Value = 4.2
$ws.Range("G15").Value = 5
$ws.Range("I15").Value = 2.12
$ws.Range("S15").Value = 2.74
$ws.Range("V15").Value = 1.89
$ws.Range("W15").Value = 1.25
$ws.Range("X15").Value = 90
$ws.Range("AC16").Value = 7.8
$ws.Range("F16").Value = 2.14
$ws.Range("G16").Value = 2.34
$ws.Range("W16").Value = 1.75
$ws.Range("F17").Value = 1.9
$ws.Range("H17").Value = 3.35
$ws.Range("I17").Value = 3.8
$ws.Range("K17").Value = 5.7
$ws.Range("S17").Value = 1.84
$ws.Range("L18").Value = 1.25
$ws.Range("I19").Value = 6
$ws.Range("V19").Value = 1.21
$ws.Range("AD20").Value = 18
$ws.Range("F20").Value = 2.44
$ws.Range("I20").Value = 3.7
$ws.Range("T20").Value = 2.08
$ws.Range("AD21").Value = 14
$ws.Range("AE21").Value = 40
$ws.Range("AF21").Value = 21
$ws.Range("AG21").Value = 13.5
$ws.Range("AJ21").Value = 50
$ws.Range("AK21").Value = 40
$ws.Range("AL21").Value = 55
$ws.Range("AN21").Value = 36
$ws.Range("AO21").Value = 40
$ws.Range("H21").Value = 2.84
$ws.Range("W21").Value = 1.56
$ws.Range("X21").Value = 12.5
$ws.Range("Y21").Value = 11
$ws.Range("Z21").Value = 21
$ws.Range("AF22").Value = 19
$ws.Range("F22").Value = 2.96
$ws.Range("G22").Value = 2.98
$ws.Range("U22").Value = 1.93
$ws.Range("W22").Value = 1.5
$ws.Range("Z22").Value = 19
$ws.Range("AH23").Value = 14.5
$ws.Range("AK23").Value = 16
$ws.Range("AM23").Value = 50
$ws.Range("F23").Value = 1.92
$ws.Range("G23").Value = 1.93
$ws.Range("V23").Value = 1.31
$ws.Range("G24").Value = 2.36
$ws.Range("H24").Value = 3.5
$ws.Range("J24").Value = 3.4
$ws.Range("W24").Value = 1.73
$ws.Range("H25").Value = 2.6
$ws.Range("W25").Value = 1.5
